$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,5,7,8,11,15,16,17,23,33,34,37)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "lipid/free"
}
